$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    22.11762717781596,
    33.9660272262802,
    43.60027209447989,
    50.54485097393651,
    55.34780330285171,
    57.55476682982366,
    58.50349776193204,
    58.35476861619483,
    57.42914312388932,
    56.06356341679849,
    54.30225190802328,
    52.48656516387938,
    50.45953388201623
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
